# Add "Luiz Felipe Brandão da Silva – RM 83025" as a new opening line of the
# document, followed by two manual line breaks, and relocate the existing
# "_GoBack" bookmark so that it sits right after the new name/RM text
# (it previously lived at the very end of the document).

$d = $word.ActiveDocument

$name = "Luiz Felipe Brandão da Silva – RM 83025"

# Step 1/2: insert two line breaks at the very start of the document (they
# land right before the existing "BrainStorm..." content). InsertBreak on a
# collapsed range always targets the start of its paragraph, so doing this
# before inserting the name text guarantees the breaks end up immediately
# before the old content.
$brAnchor = $d.Range(0, 1)
$brAnchor.Collapse(1)
$brAnchor.InsertBreak(6)

$brAnchor = $d.Range(0, 1)
$brAnchor.Collapse(1)
$brAnchor.InsertBreak(6)

# Step 3: insert the name/RM text before the two line breaks we just added.
$nameAnchor = $d.Range(0, 1)
$nameAnchor.Collapse(1)
$nameAnchor.InsertBefore($name)

# Step 4: move (or create) the "_GoBack" bookmark to sit right after the
# name text and before the two line breaks. Adding a bookmark with a name
# that already exists elsewhere in the document relocates it.
$nameLen = $name.Length
$bmRange = $d.Range($nameLen, $nameLen + 1)
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
